$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining expense values for row 29 (date 43818)
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 15
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 2
$ws.Range("M29").Value = 3

# Move the active selection to N29
$ws.Range("N29").Select()
